# Generate Report for Handoff
# Adds two new file entries (55f492e8-840b-4e7d-9ea4-c2d609985472 and
# a0cd6b83-c94f-4912-ab35-074717f3e7ca) to the Overview / zh-cn / de-de
# localization-status sheets, inserted just before the trailing
# ".localization-config" row, with matching hyperlinks.

$wb = $excel.ActiveWorkbook

$guid1 = "55f492e8-840b-4e7d-9ea4-c2d609985472"
$guid2 = "a0cd6b83-c94f-4912-ab35-074717f3e7ca"
$hash1 = "d747a0a04676c2f14efe3e265c5712bc5b03b8ee"
$hash2 = "b66b8e318a923b446cc7e0ec8f980fd5102301b0"

$md1 = "$guid1.md"
$md2 = "$guid2.md"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/d632f426acf7eb5f7f41f79b4bb3275458167373/e2e"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d632f426acf7eb5f7f41f79b4bb3275458167373/.localization-config"

function Clear-RowHyperlinks($ws, $row) {
    foreach ($hl in @($ws.Hyperlinks)) {
        if ($hl.Range.Row -eq $row) {
            $hl.Delete()
        }
    }
}

# ---------------------------------------------------------------------
# Sheet 1: "Overview"  (columns A=File Name, B=zh-cn, C=de-de)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Make room: push the trailing ".localization-config" row down two rows.
$ws.Rows("4:5").Insert()

$ws.Range("A4").Value = $md1
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

$ws.Range("A5").Value = $md2
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

# Rewire hyperlinks: rows 4-6 all need a hyperlink (4 & 5 are new, 6 is the
# shifted-down .localization-config row whose old hyperlink (still pointing
# at row 4) must be dropped and re-added at its new location).
Clear-RowHyperlinks $ws 4
Clear-RowHyperlinks $ws 5
Clear-RowHyperlinks $ws 6

$ws.Hyperlinks.Add($ws.Range("A4"), "$mdBase/$md1", "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdBase/$md2", "", "", $md2)
$ws.Hyperlinks.Add($ws.Range("A6"), $cfgUrl, "", "", ".localization-config")

$ws.Range("A4:A6").Font.Underline = $true
$ws.Range("A4:A6").Font.Color = 15570276

# ---------------------------------------------------------------------
# Sheets 2 & 3: "zh-cn" / "de-de" (the per-locale detail tables)
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Hash1Date = "2016-03-10 12:44:25"; Hash2Date = "2016-03-10 12:44:25"; OlOrg = "oltest.zh-cn" },
    @{ Name = "de-de"; Hash1Date = "2016-03-10 12:44:28"; Hash2Date = "2016-03-10 12:44:28"; OlOrg = "oltest.de-de" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Name)

    $xlf1 = "$guid1.$hash1.$($loc.Name).xlf"
    $xlf2 = "$guid2.$hash2.$($loc.Name).xlf"

    $ws.Rows("4:5").Insert()

    # New row for guid1
    $ws.Range("A4").Value = $md1
    $ws.Range("B4").Value = "Ready for handoff"
    $ws.Range("C4").Value = $xlf1
    $ws.Range("D4").Value = $loc.Hash1Date
    $ws.Range("G4").Value = "0001-01-01 00:00:00"
    $ws.Range("H4").Value = "Include"

    # New row for guid2
    $ws.Range("A5").Value = $md2
    $ws.Range("B5").Value = "Ready for handoff"
    $ws.Range("C5").Value = $xlf2
    $ws.Range("D5").Value = $loc.Hash2Date
    $ws.Range("G5").Value = "0001-01-01 00:00:00"
    $ws.Range("H5").Value = "Include"

    # Row 6 (".localization-config", pushed down from row 4) keeps its
    # original values automatically via the row insert/shift; nothing else
    # to set there.

    Clear-RowHyperlinks $ws 4
    Clear-RowHyperlinks $ws 5
    Clear-RowHyperlinks $ws 6

    $xlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash1/ol-handoff/OpenLocalizationTestOrg/$($loc.OlOrg)/ci/ht"
    $xlfBase2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash2/ol-handoff/OpenLocalizationTestOrg/$($loc.OlOrg)/ci/ht"

    $ws.Hyperlinks.Add($ws.Range("A4"), "$mdBase/$md1", "", "", $md1)
    $ws.Hyperlinks.Add($ws.Range("C4"), "$xlfBase/$xlf1", "", "", $xlf1)
    $ws.Hyperlinks.Add($ws.Range("A5"), "$mdBase/$md2", "", "", $md2)
    $ws.Hyperlinks.Add($ws.Range("C5"), "$xlfBase2/$xlf2", "", "", $xlf2)
    $ws.Hyperlinks.Add($ws.Range("A6"), $cfgUrl, "", "", ".localization-config")

    $ws.Range("A4:A6").Font.Underline = $true
    $ws.Range("A4:A6").Font.Color = 15570276
    $ws.Range("C4:C5").Font.Underline = $true
    $ws.Range("C4:C5").Font.Color = 15570276
}
